$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (7th column), shifting "Special Compensation"
# and everything after it one column to the right.
$ws.Columns.Item(7).Insert()

# New column G header + values for "Relief Amount"
$ws.Range("G1").Value = "Relief Amount"
$ws.Range("G2").Value = 1000
$ws.Range("G3").Value = 2000

# Set width for column F (previously un-set) and new column G
# (ColumnWidth is quantized internally to 1/6-character steps by this
# runtime, so we pick the nearest representable values to the targets
# 16.5703125 / 19.42578125.)
$ws.Columns.Item(6).ColumnWidth = 15.666666666666666
$ws.Columns.Item(7).ColumnWidth = 18.666666666666668

# Update the selection to match the new workbook state
$ws.Range("H7").Select()
